$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in [ 'Bluetick Hound', 'Norfolk Terrier', 'Scottish Terrier']and diag.disease_term in ['Bladder Cancer'] and demo.sex in ['Female'] and demo.neutered_indicator IN ['Yes']
    
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$statsQuery = @'
MATCH (s:study)
  MATCH (demo:demographic) 
  MATCH (diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in [ 'Bluetick Hound', 'Norfolk Terrier', 'Scottish Terrier']and diag.disease_term in ['Bladder Cancer'] and demo.sex in ['Female'] and demo.neutered_indicator IN ['Yes']
    
OPTIONAL MATCH (s)<-[:member_of]-(c:case)
OPTIONAL MATCH (c)<-[:of_case]-(samp:sample)<-[:of_sample]-(f:file)
RETURN 
	count(DISTINCT(f)) as number_of_files , 
	count(DISTINCT(samp)) as number_of_sample , 
	count(DISTINCT(c.case_id)) as number_of_cases , 
	count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@

$neo4jName = @'
TC21_Canine_StudyUBC02-Breed_Sex_NeuteredStatus_Neo4jData.xlsx
'@

$webDataName = @'
TC21_Canine_StudyUBC02-Breed_Sex_NeuteredStatus_WebData.xlsx
'@

$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statsQuery
$ws.Range("D2").Value = $neo4jName
$ws.Range("E2").Value = $webDataName

$ws.Rows.Item(2).RowHeight = 304.5
$ws.Rows.Item(3).RowHeight = 275.5
$ws.Rows.Item(4).RowHeight = 275.5

$ws.Range("B2").Select()
